# Tetris board edit: move a falling piece's landing spots and clear a few
# bottom-stack cells (per "excel to test rotation and path").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("visual")

# Reference cells carrying the styles we need to reuse (so we don't create
# brand-new style entries in the workbook):
#   D2  -> style "1" (plain/blank board cell)
#   AG3 -> style "7" (gray filled "." block)
#   A4  -> style "8" (red filled "." block)

$xlPasteFormats = -4122

# --- Add new gray "." blocks (style 7) -------------------------------------
$grayTargets = "F3", "F4", "L4", "F5", "L5", "L6"
$ws.Range("AG3").Copy() | Out-Null
foreach ($addr in $grayTargets) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value = "."
}

# --- Add new red "." blocks (style 8) --------------------------------------
$redTargets = "G2", "H2", "I2", "J2"
$ws.Range("A4").Copy() | Out-Null
foreach ($addr in $redTargets) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value = "."
}

# --- Clear existing blocks back to the plain style (style 1) ---------------
$clearTargets = "F15", "K16", "I19", "E21"
$ws.Range("D2").Copy() | Out-Null
foreach ($addr in $clearTargets) {
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).ClearContents()
}

$excel.CutCopyMode = 0

# --- Update selection / scroll position ------------------------------------
$ws.Range("J2").Select()
